# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1 = "Wins", AE1 = "Losses", AF1 = "Ties" ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting (bold, border, centered) from the existing header
# cell AC1 so the new header cells match the rest of the header row,
# reusing the same style instead of creating a new one.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-53): same W/L/T values for every player on the roster ---
$lastRow = 53
$wins = $ws.Range("AD2:AD" + $lastRow)
$losses = $ws.Range("AE2:AE" + $lastRow)
$ties = $ws.Range("AF2:AF" + $lastRow)

$wins.Value = 59
$losses.Value = 103
$ties.Value = 0

$excel.CutCopyMode = 0
